{"js": "const body = context.document.body;\n\n{\n  const results = body.search(\"479\u00d77=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"673\u00d74=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"187\u00d77=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"551\u00d72=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"294\u00d75=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"682\u00d74=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"754\u00d79=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"680\u00d73=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"928\u00d72=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"271\u00d74=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"897\u00d77=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"246\u00d79=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"516\u00d78=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"471\u00d72=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"977\u00d75=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"288\u00d73=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"137\u00d76=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"559\u00d74=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"677\u00d75=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"503\u00d72=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"152\u00d74=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"832\u00d72=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"525\u00d76=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"927\u00d78=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"305\u00d74=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"850\u00d73=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"687\u00d78=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"564\u00d72=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"683\u00d74=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"282\u00d75=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"864\u00d76=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"848\u00d79=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"714\u00d74=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"462\u00d73=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"907\u00d76=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"615\u00d76=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"891\u00d78=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"541\u00d77=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"469\u00d72=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"837\u00d75=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"875\u00d76=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"358\u00d73=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"852\u00d74=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"571\u00d79=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"916\u00d77=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"564\u00d75=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"885\u00d74=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"817\u00d73=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"731\u00d75=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"827\u00d73=\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"479\u00d77=\"\n$find.Replacement.Text = \"673\u00d74=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"187\u00d77=\"\n$find.Replacement.Text = \"551\u00d72=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"294\u00d75=\"\n$find.Replacement.Text = \"682\u00d74=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"754\u00d79=\"\n$find.Replacement.Text = \"680\u00d73=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"928\u00d72=\"\n$find.Replacement.Text = \"271\u00d74=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"897\u00d77=\"\n$find.Replacement.Text = \"246\u00d79=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"516\u00d78=\"\n$find.Replacement.Text = \"471\u00d72=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"977\u00d75=\"\n$find.Replacement.Text = \"288\u00d73=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"137\u00d76=\"\n$find.Replacement.Text = \"559\u00d74=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"677\u00d75=\"\n$find.Replacement.Text = \"503\u00d72=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"152\u00d74=\"\n$find.Replacement.Text = \"832\u00d72=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"525\u00d76=\"\n$find.Replacement.Text = \"927\u00d78=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"305\u00d74=\"\n$find.Replacement.Text = \"850\u00d73=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"687\u00d78=\"\n$find.Replacement.Text = \"564\u00d72=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"683\u00d74=\"\n$find.Replacement.Text = \"282\u00d75=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"864\u00d76=\"\n$find.Replacement.Text = \"848\u00d79=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"714\u00d74=\"\n$find.Replacement.Text = \"462\u00d73=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"907\u00d76=\"\n$find.Replacement.Text = \"615\u00d76=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"891\u00d78=\"\n$find.Replacement.Text = \"541\u00d77=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"469\u00d72=\"\n$find.Replacement.Text = \"837\u00d75=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"875\u00d76=\"\n$find.Replacement.Text = \"358\u00d73=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"852\u00d74=\"\n$find.Replacement.Text = \"571\u00d79=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"916\u00d77=\"\n$find.Replacement.Text = \"564\u00d75=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"885\u00d74=\"\n$find.Replacement.Text = \"817\u00d73=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"731\u00d75=\"\n$find.Replacement.Text = \"827\u00d73=\"\n$find.MatchWildcards = $false\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n"}
